# Fruta / hortaliza, semanal
# Weekly update: a new price observation (week of 2022-11-29, serial 44894)
# is inserted as the new row 16 for "Feria Lagunitas de Puerto Montt -
# Arándano (blue)", pushing the existing rows 16:23 down to 17:24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 16; rows 16:23 shift down to 17:24
# (dimension grows from A1:T23 to A1:T24 automatically).
$ws.Rows("16:16").Insert()

# Populate the new row 16 with this week's data.
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C16").Value = 'Los Lagos'
$ws.Range("D16").Value = 44894
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100101
$ws.Range("H16").Value = 'Berries'
$ws.Range("I16").Value = 100101001
$ws.Range("J16").Value = 'Arándano (blue)'
$ws.Range("K16").Value = 'Sin especificar'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8500
$ws.Range("P16").Value = 8250
$ws.Range("Q16").Value = '$/bandeja 2 kilos'
$ws.Range("R16").Value = 'Provincia de Curicó'
$ws.Range("S16").Value = 4125
$ws.Range("T16").Value = 2
